# feat/CIV-943 sot format updated
#
# The statement-of-truth placeholder "<<statementOfTruth.name>>" (in the
# "Statement of truth" table) needs to be rendered in the GDSTransportWebsite
# font, matching the rest of the surrounding body text.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("<<statementOfTruth.name>>", $false, $false, $false,
                            $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Font.Name = "GDSTransportWebsite"
}
